$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 10:28:12'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 171'
$ws1.Cells.Item(110, 1).Value = '07:59:05'
$ws1.Cells.Item(110, 2).Value = '09:23'
$ws1.Cells.Item(110, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(110, 4).Value = 84
$ws1.Cells.Item(110, 5).Value = 'LP1912'
$ws1.Cells.Item(111, 1).Value = '07:46:15'
$ws1.Cells.Item(111, 2).Value = '09:23'
$ws1.Cells.Item(111, 3).Value = '17_ROMERO'
$ws1.Cells.Item(111, 4).Value = 97
$ws1.Cells.Item(111, 5).Value = 'LP1912'
$ws1.Cells.Item(112, 1).Value = '08:21:27'
$ws1.Cells.Item(112, 2).Value = '09:23'
$ws1.Cells.Item(112, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(112, 4).Value = 62
$ws1.Cells.Item(112, 5).Value = 'LP1912'
$ws1.Cells.Item(120, 1).Value = '08:50:00'
$ws1.Cells.Item(120, 2).Value = '09:35'
$ws1.Cells.Item(120, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(120, 4).Value = 45
$ws1.Cells.Item(120, 5).Value = 'LP1912'
$ws1.Cells.Item(121, 1).Value = '08:57:11'
$ws1.Cells.Item(121, 2).Value = '09:35'
$ws1.Cells.Item(121, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(121, 4).Value = 38
$ws1.Cells.Item(121, 5).Value = 'LP1912'
$ws1.Cells.Item(144, 1).Value = '10:28:12'
$ws1.Cells.Item(144, 2).Value = '10:29'
$ws1.Cells.Item(144, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(144, 4).Value = 1
$ws1.Cells.Item(144, 5).Value = 'LP1912'
$ws1.Cells.Item(145, 1).Value = '10:28:12'
$ws1.Cells.Item(145, 2).Value = '10:30'
$ws1.Cells.Item(145, 3).Value = '10_OLMOS'
$ws1.Cells.Item(145, 4).Value = 2
$ws1.Cells.Item(145, 5).Value = 'LP1912'
$ws1.Cells.Item(146, 1).Value = '10:28:12'
$ws1.Cells.Item(146, 2).Value = '10:35'
$ws1.Cells.Item(146, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(146, 4).Value = 7
$ws1.Cells.Item(146, 5).Value = 'LP1912'
$ws1.Cells.Item(147, 1).Value = '08:50:00'
$ws1.Cells.Item(147, 2).Value = '10:42'
$ws1.Cells.Item(147, 3).Value = '17_ROMERO'
$ws1.Cells.Item(147, 4).Value = 112
$ws1.Cells.Item(147, 5).Value = 'LP1912'
$ws1.Cells.Item(148, 1).Value = '08:50:00'
$ws1.Cells.Item(148, 2).Value = '10:44'
$ws1.Cells.Item(148, 3).Value = '14_ABASTO'
$ws1.Cells.Item(148, 4).Value = 114
$ws1.Cells.Item(148, 5).Value = 'LP1912'
$ws1.Cells.Item(149, 1).Value = '09:38:04'
$ws1.Cells.Item(149, 2).Value = '10:52'
$ws1.Cells.Item(149, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(149, 4).Value = 74
$ws1.Cells.Item(149, 5).Value = 'LP1912'
$ws1.Cells.Item(150, 1).Value = '10:28:12'
$ws1.Cells.Item(150, 2).Value = '10:52'
$ws1.Cells.Item(150, 3).Value = '15_ABASTO'
$ws1.Cells.Item(150, 4).Value = 24
$ws1.Cells.Item(150, 5).Value = 'LP1912'
$ws1.Cells.Item(151, 1).Value = '10:28:12'
$ws1.Cells.Item(151, 2).Value = '10:53'
$ws1.Cells.Item(151, 3).Value = '10_OLMOS'
$ws1.Cells.Item(151, 4).Value = 25
$ws1.Cells.Item(151, 5).Value = 'LP1912'
$ws1.Cells.Item(152, 1).Value = '10:28:12'
$ws1.Cells.Item(152, 2).Value = '10:57'
$ws1.Cells.Item(152, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(152, 4).Value = 29
$ws1.Cells.Item(152, 5).Value = 'LP1912'
$ws1.Cells.Item(153, 1).Value = '10:28:12'
$ws1.Cells.Item(153, 2).Value = '10:57'
$ws1.Cells.Item(153, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(153, 4).Value = 29
$ws1.Cells.Item(153, 5).Value = 'LP1912'
$ws1.Cells.Item(154, 1).Value = '09:38:04'
$ws1.Cells.Item(154, 2).Value = '11:02'
$ws1.Cells.Item(154, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(154, 4).Value = 84
$ws1.Cells.Item(154, 5).Value = 'LP1912'
$ws1.Cells.Item(155, 1).Value = '10:28:12'
$ws1.Cells.Item(155, 2).Value = '11:04'
$ws1.Cells.Item(155, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(155, 4).Value = 36
$ws1.Cells.Item(155, 5).Value = 'LP1912'
$ws1.Cells.Item(156, 1).Value = '10:28:12'
$ws1.Cells.Item(156, 2).Value = '11:05'
$ws1.Cells.Item(156, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(156, 4).Value = 37
$ws1.Cells.Item(156, 5).Value = 'LP1912'
$ws1.Cells.Item(157, 1).Value = '09:38:04'
$ws1.Cells.Item(157, 2).Value = '11:07'
$ws1.Cells.Item(157, 3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(157, 4).Value = 89
$ws1.Cells.Item(157, 5).Value = 'LP1912'
$ws1.Cells.Item(158, 1).Value = '09:38:04'
$ws1.Cells.Item(158, 2).Value = '11:12'
$ws1.Cells.Item(158, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(158, 4).Value = 94
$ws1.Cells.Item(158, 5).Value = 'LP1912'
$ws1.Cells.Item(159, 1).Value = '10:28:12'
$ws1.Cells.Item(159, 2).Value = '11:12'
$ws1.Cells.Item(159, 3).Value = '15_ABASTO'
$ws1.Cells.Item(159, 4).Value = 44
$ws1.Cells.Item(159, 5).Value = 'LP1912'
$ws1.Cells.Item(160, 1).Value = '09:38:04'
$ws1.Cells.Item(160, 2).Value = '11:20'
$ws1.Cells.Item(160, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws1.Cells.Item(160, 4).Value = 102
$ws1.Cells.Item(160, 5).Value = 'LP1912'
$ws1.Cells.Item(161, 1).Value = '10:28:12'
$ws1.Cells.Item(161, 2).Value = '11:21'
$ws1.Cells.Item(161, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(161, 4).Value = 53
$ws1.Cells.Item(161, 5).Value = 'LP1912'
$ws1.Cells.Item(162, 1).Value = '09:38:04'
$ws1.Cells.Item(162, 2).Value = '11:25'
$ws1.Cells.Item(162, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(162, 4).Value = 107
$ws1.Cells.Item(162, 5).Value = 'LP1912'
$ws1.Cells.Item(163, 1).Value = '09:38:04'
$ws1.Cells.Item(163, 2).Value = '11:27'
$ws1.Cells.Item(163, 3).Value = '225_C ROCA-H SUR'
$ws1.Cells.Item(163, 4).Value = 109
$ws1.Cells.Item(163, 5).Value = 'LP1912'
$ws1.Cells.Item(164, 1).Value = '09:38:04'
$ws1.Cells.Item(164, 2).Value = '11:32'
$ws1.Cells.Item(164, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(164, 4).Value = 114
$ws1.Cells.Item(164, 5).Value = 'LP1912'
$ws1.Cells.Item(165, 1).Value = '09:38:04'
$ws1.Cells.Item(165, 2).Value = '11:36'
$ws1.Cells.Item(165, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(165, 4).Value = 118
$ws1.Cells.Item(165, 5).Value = 'LP1912'
$ws1.Cells.Item(166, 1).Value = '10:28:12'
$ws1.Cells.Item(166, 2).Value = '11:40'
$ws1.Cells.Item(166, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(166, 4).Value = 72
$ws1.Cells.Item(166, 5).Value = 'LP1912'
$ws1.Cells.Item(167, 1).Value = '10:28:12'
$ws1.Cells.Item(167, 2).Value = '11:42'
$ws1.Cells.Item(167, 3).Value = '17_ROMERO'
$ws1.Cells.Item(167, 4).Value = 74
$ws1.Cells.Item(167, 5).Value = 'LP1912'
$ws1.Cells.Item(168, 1).Value = '10:28:12'
$ws1.Cells.Item(168, 2).Value = '11:51'
$ws1.Cells.Item(168, 3).Value = '10_OLMOS'
$ws1.Cells.Item(168, 4).Value = 83
$ws1.Cells.Item(168, 5).Value = 'LP1912'
$ws1.Cells.Item(169, 1).Value = '10:28:12'
$ws1.Cells.Item(169, 2).Value = '11:51'
$ws1.Cells.Item(169, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(169, 4).Value = 83
$ws1.Cells.Item(169, 5).Value = 'LP1912'
$ws1.Cells.Item(170, 1).Value = '10:28:12'
$ws1.Cells.Item(170, 2).Value = '11:59'
$ws1.Cells.Item(170, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(170, 4).Value = 91
$ws1.Cells.Item(170, 5).Value = 'LP1912'
$ws1.Cells.Item(171, 1).Value = '10:28:12'
$ws1.Cells.Item(171, 2).Value = '12:06'
$ws1.Cells.Item(171, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(171, 4).Value = 98
$ws1.Cells.Item(171, 5).Value = 'LP1912'
$ws1.Cells.Item(172, 1).Value = '10:28:12'
$ws1.Cells.Item(172, 2).Value = '12:07'
$ws1.Cells.Item(172, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(172, 4).Value = 99
$ws1.Cells.Item(172, 5).Value = 'LP1912'
$ws1.Cells.Item(173, 1).Value = '10:28:12'
$ws1.Cells.Item(173, 2).Value = '12:14'
$ws1.Cells.Item(173, 3).Value = '17_ROMERO'
$ws1.Cells.Item(173, 4).Value = 106
$ws1.Cells.Item(173, 5).Value = 'LP1912'
$ws1.Cells.Item(174, 1).Value = '10:28:12'
$ws1.Cells.Item(174, 2).Value = '12:18'
$ws1.Cells.Item(174, 3).Value = '14_ABASTO'
$ws1.Cells.Item(174, 4).Value = 110
$ws1.Cells.Item(174, 5).Value = 'LP1912'
$ws1.Cells.Item(175, 1).Value = '10:28:12'
$ws1.Cells.Item(175, 2).Value = '12:21'
$ws1.Cells.Item(175, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(175, 4).Value = 113
$ws1.Cells.Item(175, 5).Value = 'LP1912'
$ws1.Cells.Item(176, 1).Value = '10:28:12'
$ws1.Cells.Item(176, 2).Value = '12:21'
$ws1.Cells.Item(176, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(176, 4).Value = 113
$ws1.Cells.Item(176, 5).Value = 'LP1912'

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 10:28:12'
$ws2.Cells.Item(3, 1).Value = 'Total filas: 22'
$ws2.Cells.Item(26, 1).Value = '10:28:12'
$ws2.Cells.Item(26, 2).Value = '11:51'
$ws2.Cells.Item(26, 3).Value = '215B_EL PATO'
$ws2.Cells.Item(26, 4).Value = 83
$ws2.Cells.Item(26, 5).Value = 'LP1912'
$ws2.Cells.Item(27, 1).Value = '10:28:12'
$ws2.Cells.Item(27, 2).Value = '12:21'
$ws2.Cells.Item(27, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(27, 4).Value = 113
$ws2.Cells.Item(27, 5).Value = 'LP1912'

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 10:28:12'
$ws3.Cells.Item(3, 1).Value = 'Total filas: 31'
$ws3.Cells.Item(36, 1).Value = '10:28:12'
$ws3.Cells.Item(36, 2).Value = '12:04'
$ws3.Cells.Item(36, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(36, 4).Value = 96
$ws3.Cells.Item(36, 5).Value = 'L6173'
